$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.047.44'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.159.79'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +2.40%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.24'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.13'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.86%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.157.51'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.528'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.80%  '
$ws.Range('E10').Value = '  +3.16%  '
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('E12').Value = '  +2.18%  '
$ws.Range('E13').Value = '  +13.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.18'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.676.65'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.082.92'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.157.70'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('E18').Value = '  +3.32%  '
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '506.21'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.77%  '
$ws.Range('E21').Value = '  +2.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.719'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.32'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.58%  '
$ws.Range('E24').Value = '  +1.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.39'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.91'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +6.93%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.90'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.53%  '
$ws.Range('E29').Value = '  +3.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '27.62'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.01%  '
$ws.Range('E31').Value = '  +6.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('E33').Value = '  +1.97%  '
$ws.Range('E34').Value = '  +5.22%  '
$ws.Range('E35').Value = '  +2.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '54.92'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.65%  '
$ws.Range('E37').Value = '  +8.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '466.60'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0419'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.99'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.85%  '
$ws.Range('E41').Value = '  +3.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.046.68'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.45'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +7.29%  '
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.53'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.48%  '
$ws.Range('E47').Value = '  +11.15%  '
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.115'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('E50').Value = '  +3.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.48'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.11%  '
